$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix the two "Free before X/after Y" availability entries to the
# (bug-fixed) "Free except X-Y" phrasing.
$ws.Range("D5").Value = "Free except 4-5:30"
$ws.Range("G5").Value = "Free except 3:30-4:30"

# Remove Brandon's row entirely (row 7) - everything below shifts up.
$ws.Range("A7:I7").EntireRow.Delete()

# Select the row that now holds the next entry (Johanne), matching the
# saved selection state in the workbook.
$ws.Rows.Item(7).Select()

# Columns C and H grew slightly wider once the data reflowed (column C now
# also carries the longer "Free after 5:30pm" text, and H carries "before
# 9:30pm"); match the resulting best-fit widths.
$ws.Columns.Item(3).ColumnWidth = 14.1666666667
$ws.Columns.Item(8).ColumnWidth = 12.1666666667
